$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.011.56"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "1.670.89"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("D12").Value = "1.906.37"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "1.666.39"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.521"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "27.031.98"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "0.0₃0736"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "1.452.24"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.572"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.894"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.29%  "
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +11.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.79%  "
$ws.Range("D45").Value = "1.815.15"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.52%  "
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  +1.65%  "
